# Updated cryptos list (price/volume refresh + PancakeSwap/PEPE row swap).
#
# Price (column D) and Volume(1h) (column E) are stored as literal text in
# this sheet (e.g. "163.00", "66.173.33", "  +5.97%  "), not as numbers, so
# trailing zeros / multi-dot "thousands" groupings / padding survive. Any D
# value that Excel would otherwise auto-parse into a number is written with
# NumberFormat "@" (Text) set first so it round-trips as the exact string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.173.33"
$ws.Range("E2").Value = "  +5.97%  "

# Row 3
$ws.Range("D3").Value = "2.997.93"
$ws.Range("E3").Value = "  +3.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.51"
$ws.Range("E5").Value = "  +2.60%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "163.00"
$ws.Range("E6").Value = "  +13.43%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  +3.57%  "

# Row 9
$ws.Range("D9").Value = "2.995.60"
$ws.Range("E9").Value = "  +3.31%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.58"
$ws.Range("E10").Value = "  -4.76%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.155"
$ws.Range("E11").Value = "  +4.07%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.456"
$ws.Range("E12").Value = "  +5.71%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000252"
$ws.Range("E13").Value = "  +6.44%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.57"
$ws.Range("E14").Value = "  +5.71%  "

# Row 15
$ws.Range("E15").Value = "  -0.69%  "

# Row 16
$ws.Range("D16").Value = "66.130.82"
$ws.Range("E16").Value = "  +5.82%  "

# Row 17
$ws.Range("D17").Value = "3.492.56"
$ws.Range("E17").Value = "  +3.21%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.91"
$ws.Range("E18").Value = "  +4.86%  "

# Row 19
$ws.Range("D19").Value = "2.995.55"
$ws.Range("E19").Value = "  +3.35%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "453.15"
$ws.Range("E20").Value = "  +6.47%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.84"
$ws.Range("E21").Value = "  +6.26%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.686"
$ws.Range("E22").Value = "  +4.43%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.35"
$ws.Range("E23").Value = "  +7.24%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "82.33"
$ws.Range("E24").Value = "  +4.72%  "

# Row 25
$ws.Range("E25").Value = "  +13.90%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.26"
$ws.Range("E26").Value = "  +3.43%  "

# Row 27
$ws.Range("E27").Value = "  +3.49%  "

# Row 28
$ws.Range("E28").Value = "  +0.03%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.10"
$ws.Range("E29").Value = "  +14.72%  "

# Row 30
$ws.Range("E30").Value = "  +19.25%  "

# Row 31
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0000104"
$ws.Range("E31").Value = "  -4.64%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.62"
$ws.Range("E32").Value = "  +5.66%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "27.30"
$ws.Range("E33").Value = "  +6.17%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.111"
$ws.Range("E34").Value = "  +5.17%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.999"
$ws.Range("E35").Value = "  -0.29%  "

# Row 36
$ws.Range("E36").Value = "  +4.44%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.81"
$ws.Range("E37").Value = "  +8.42%  "

# Row 38
$ws.Range("E38").Value = "  +8.63%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "49.56"
$ws.Range("E39").Value = "  +2.15%  "

# Row 40
$ws.Range("E40").Value = "  +1.63%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.310"
$ws.Range("E41").Value = "  +16.80%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "44.35"
$ws.Range("E42").Value = "  +7.51%  "

# Row 43
$ws.Range("E43").Value = "  +6.91%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.47"
$ws.Range("E44").Value = "  +5.42%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "399.24"
$ws.Range("E45").Value = "  +12.33%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0359"
$ws.Range("E46").Value = "  +6.43%  "

# Row 47
$ws.Range("D47").Value = "2.770.75"
$ws.Range("E47").Value = "  +2.15%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "133.27"
$ws.Range("E48").Value = "  +0.52%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "23.85"
$ws.Range("E50").Value = "  +12.40%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.107"
$ws.Range("E51").Value = "  +4.05%  "
